$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 107 (shifts existing rows 107..205 down to 108..206)
$ws.Rows(107).Insert()

# Populate the newly inserted row 107 with the new record's data
$ws.Range("A107").Value = 9
$ws.Range("B107").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C107").Value = "Metropolitana"
$ws.Range("D107").Value = 44680
$ws.Range("E107").Value = 13
$ws.Range("F107").Value = 100112026
$ws.Range("G107").Value = "Haba"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 52
$ws.Range("K107").Value = 20000
$ws.Range("L107").Value = 21000
$ws.Range("M107").Value = 20500
$ws.Range("N107").Value = "`$/saco 25 kilos"
$ws.Range("O107").Value = "Provincia del Elquí"
$ws.Range("P107").Value = 820
$ws.Range("Q107").Value = 25
$ws.Range("R107").Value = "Hortaliza"
